$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit permutes the per-record data of rows 2-15 (each row keeps its
# constant/shared columns - C, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG,
# AW, AX - in place) while the record-specific columns (A, B, D, E, F, G, H,
# Q, R) and the single "Gamla hack" remark (column AC) move to different
# rows, as if the underlying records were reordered.
#
# Mapping below: for each destination row, which row currently (before the
# edit) holds the content that should end up there.
$sourceRowFor = @{
    2  = 10
    3  = 5
    4  = 15
    5  = 11
    6  = 8
    7  = 12
    8  = 2
    9  = 7
    10 = 6
    11 = 13
    12 = 4
    13 = 14
    14 = 3
    15 = 9
}

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# 1) Snapshot all current values first (reads must happen before any writes,
#    since several destinations read from rows that are themselves about to
#    be overwritten).
$snapshot = @{}
foreach ($col in $cols) {
    $rowVals = @{}
    for ($r = 2; $r -le 15; $r++) {
        $rowVals[$r] = $ws.Range($col + $r).Value2
    }
    $snapshot[$col] = $rowVals
}
$acSnapshot = @{}
for ($r = 2; $r -le 15; $r++) {
    $acSnapshot[$r] = $ws.Range("AC" + $r).Value2
}

# 2) Write the permuted values back out.
foreach ($destRow in $sourceRowFor.Keys) {
    $srcRow = $sourceRowFor[$destRow]
    foreach ($col in $cols) {
        $ws.Range($col + $destRow).Value2 = $snapshot[$col][$srcRow]
    }
    $acValue = $acSnapshot[$srcRow]
    if ($acValue -ne $null -and $acValue -ne "") {
        $ws.Range("AC" + $destRow).Value2 = $acValue
    } else {
        $ws.Range("AC" + $destRow).ClearContents()
    }
}
